$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns
# D-column values that look like plain numbers must be forced to text
# (NumberFormat '@' then reset the style back to Normal so no stray
# style index is left attached to the cell) so they keep their exact
# textual representation (leading/trailing zeros, decimal grouping).

$ws.Range('D2').Value = '28.953.23'
$ws.Range('E2').Value = '  +2.02%  '
$ws.Range('D3').Value = '1.904.62'
$ws.Range('E3').Value = '  +2.02%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '332.65'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4642'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4051'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.35%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.96'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07995'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.003'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('E12').Value = '  -0.87%  '
$ws.Range('D13').Value = '1.909.79'
$ws.Range('E13').Value = '  +2.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.921'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.056'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.25%  '
$ws.Range('E16').Value = '  -0.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.70'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.68%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001033'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.73%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06572'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.61%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.42'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.001'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').Value = '28.975.83'
$ws.Range('E22').Value = '  +2.06%  '
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('E24').Value = '  +1.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.237'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.44%  '
$ws.Range('D26').Value = '2.131.98'
$ws.Range('E26').Value = '  +2.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '157.67'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.71'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.097'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.26%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.390'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '118.73'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9790'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.32%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09385'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.416'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.40%  '
$ws.Range('E36').Value = '  -1.04%  '
$ws.Range('E37').Value = '  -0.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02223'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.76%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.383'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.27%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.164'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.00%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.001'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5786'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.36%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1822'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.37%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.13'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.89%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.262'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.329'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +13.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.08'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5491'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.34%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.898'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.77%  '
$ws.Range('E50').Value = '  +2.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '47.86'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +23.82%  '
